$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 11 new rows above row 74 (old rows 74-81 shift down to 85-92).
# Excel's row-insert copies formatting from the row above, which matches
# the target styling (A/B -> style 1, C:M -> style 3) already used by row 73.
$ws.Range("A74:A84").EntireRow.Insert()

# New "late timepoint" CAR genotype rows (days 17, 22, 26) inserted at 74-84.
$newRows = @(
    @{B=17; C=102;   D=47;   E=4.8899999999999997;  F=73.599999999999994; G=2.35;                H=5.95; I=2.7;  J=3.56; K=89.6;                L=0.35; M=14.6},
    @{B=17; C=85;    D=53;   E=3.81;                 F=74.2;               G=1.67;                H=4.92; I=4.54; J=3.85; K=65.400000000000006;  L=0.4;  M=17.8},
    @{B=17; C=83.7;  D=54.1; E=4.12;                 F=71.8;               G=3.05;                H=5.28; I=4.49; J=2.98; K=70.3;                 L=0.42; M=13.1},
    @{B=22; C=75;    D=53.3; E=2.82;                 F=73.599999999999994; G=1.66;                H=4.55; I=2.25; J=1.88; K=75.599999999999994;  L=0.37; M=14.1},
    @{B=22; C=91.5;  D=51.6; E=3;                    F=74.900000000000006; G=1.68;                H=5.74; I=2.57; J=2;    K=82.8;                 L=0.47; M=14.6},
    @{B=22; C=86.5;  D=50.8; E=2.2000000000000002;   F=73.3;               G=2.29;                H=5.05; I=5.17; J=1.49; K=71.7;                 L=0.36; M=13.4},
    @{B=22; C=87;    D=56.7; E=4.0999999999999996;   F=75.099999999999994; G=2.71;                H=6.11; I=4.12; J=2.6;  K=89.1;                 L=0.46; M=18.899999999999999},
    @{B=26; C=59;    D=56.1; E=3.3;                  F=70.2;               G=2.31;                H=6.57; I=3.06; J=2.02; K=78.8;                 L=0.36; M=16.7},
    @{B=26; C=60.5;  D=53.8; E=3.8;                  F=68;                 G=2.4900000000000002;  H=7.1;  I=5.0199999999999996; J=1.71; K=92.3;    L=0.3;  M=21.2},
    @{B=26; C=65.5;  D=54.3; E=3.8;                  F=69.8;               G=2.5;                 H=7.84; I=5.52; J=1.7;  K=84.4;                 L=0.35; M=16.899999999999999},
    @{B=26; C=70;    D=53.7; E=2.6;                  F=69.8;               G=1.8;                 H=6.14; I=2.4500000000000002; J=1.75; K=68.5;    L=0.43; M=19.2}
)

$r = 74
foreach ($row in $newRows) {
    $ws.Range("A$r").Value = "CAR"
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
    $ws.Range("H$r").Value = $row.H
    $ws.Range("I$r").Value = $row.I
    $ws.Range("J$r").Value = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
    $ws.Range("M$r").Value = $row.M
    $r = $r + 1
}

# Sheet view / selection updates captured in the diff.
$ws.Application.ActiveWindow.Zoom = 188
$ws.Range("A1:M1").Select()

Write-Output "done"
